$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer car model names.
# NOTE: the host engine stores the OOXML <col width> as ColumnWidth + 5/6
# (its emulation of Excel's character-padding conversion), so we dial the
# COM-visible ColumnWidth back by that fixed offset to land on exactly the
# target stored width of 23.
$ws.Columns.Item(1).ColumnWidth = 23 - (5/6)

# Insert a new row at position 12 (this shifts the old Toyota/Honda/NIO rows
# down to 13/14/15, matching the diff's reshuffle)
$ws.Rows.Item(12).Insert()

# Update the car-name labels in column A to include the model name
$ws.Range("A2").Value  = "Ford Tourneo Custom"
$ws.Range("A3").Value  = "ZEEKR 001"
$ws.Range("A4").Value  = "MAXUS MIFA 7"
$ws.Range("A5").Value  = "VW Passat"
$ws.Range("A6").Value  = "Škoda Kodiaq"
$ws.Range("A7").Value  = "BMW X2"
$ws.Range("A8").Value  = "Renault Rafale HEV"
$ws.Range("A9").Value  = "Mercedes-Benz E-Class"
$ws.Range("A10").Value = "Suzuki Swift"
$ws.Range("A11").Value = "Dacia Duster"

# Populate the newly inserted row 12 with the "Renault Espace" entry
# (numerically identical to the original Renault row's data)
$ws.Range("A12").Value = "Renault Espace"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 19.05
$ws.Range("D12").Value = 25.6
$ws.Range("E12").Value = 30.95
$ws.Range("F12").Value = 10.12
$ws.Range("G12").Value = 10.71
$ws.Range("H12").Value = 3.57
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 100

# Rename the shifted-down rows to their new model names
$ws.Range("A13").Value = "Toyota C-HR"
$ws.Range("A14").Value = "Honda CR-V"
$ws.Range("A15").Value = "NIO EL6"
